# Auto-generated Word COM-interop script.
# Rebuilds the CSS Notes body from the single 'CSS Notes' paragraph,
# matching the target OOXML diff.

$d = $word.ActiveDocument

# --- Append the remaining paragraphs first (plain text; formatting applied after) ---
# Inserted *before* styling paragraph 1 so the new runs don't inherit
# the title's bold/underline from the insertion point's current formatting.
$insertPoint = $d.Paragraphs(1).Range
$body = "CSS is a powerful tool that changes the presentation of a single document or multiple documents at one go. CSS was first proposed in 1994 just when web was evolving. CSS provides a simple and declarative style language that was flexible enough for the developers helping them to provide styles of the documents they are authoring. These documents can be HTML or XML or any document that was written with some type of markup language." + "`r`n" + "CSS stands for Cascading style sheets. Here cascade means authors can define multiple styles in many different ways and these styles can be cascaded to many different levels where those levels may be of different documents or different elements in the same document. By 1996, CSS1 was finished while CSS working group started working on CSS2. As CSS started evolving, many different browsers started adopting it and now the implementation of CSS standards was being provided by almost all the browsers and CSS working group is progressing towards providing more mature effects from CSS." + "`r`n" + "Elements: elements are the basis for document structure. In HTML, the most common elements are p, table, span, a, and div. Each of these elements plays a role in the display of the document" + "`r`n" + "Replaced and Non-Replaced Elements: Even though CSS depends on the elements of a document to decide how that element needs to be displayed on the screen, all the elements are not created equally. In CSS, elements take two forms" + "`r`n" + "Replaced elements: Replaced elements are those elements whose content is replaced with something else that is not directly represented by document content. For example, an image element is replaced with an image that is specified by the address of the image. Here the address of the image is the original document content but this content is replaced with the original image that is represented by that document content. Hence this is considered to be a replaced element. Some more examples of replaced elements are input elements like radio button, check box, submit button, reset button, text box, text area etc.," + "`r`n" + "Non-Replaced elements: the majority of HTML elements are non-replaced elements. That means their content is presented by the user agent. Generally the user agent will be the browser. User agent presents the content inside a box that is generated by the element itself. For example, <span>hi</span> is a non-replaced element and the text “hi” will be displayed by the user agent. This is true for paragraphs, headings, table cells, lists and almost everything else in HTML." + "`r`n" + "Element Display Roles: "
$insertPoint.InsertAfter("`r`n" + $body)

Write-Host "Paragraph count:" $d.Paragraphs.Count

# --- Paragraph 1: bold + single-underline the title (applied last so it doesn't leak) ---
$titlePara = $d.Paragraphs(1).Range
$titlePara.Font.Bold = 1
$titlePara.Font.Underline = 1

# --- Paragraph 2: split into its 2 source runs ---
$search2 = $d.Paragraphs(2).Range
$seg2_0 = $search2
$null = $seg2_0.Find.Execute("CSS is a powerful tool that changes the presentation of a single document or multiple documents at one go. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg2_0.Font.Bold = 1
$seg2_0.Font.Bold = 0
$search2 = $d.Paragraphs(2).Range
$search2.SetRange($seg2_0.End, $search2.End)
$seg2_1 = $search2
$null = $seg2_1.Find.Execute("CSS was first proposed in 1994 just when web was evolving. CSS provides a simple and declarative style language that was flexible enough for the developers helping them to provide styles of the documents they are authoring. These documents can be HTML or XML or any document that was written with some type of markup language.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg2_1.Font.Bold = 1
$seg2_1.Font.Bold = 0
$search2 = $d.Paragraphs(2).Range
$search2.SetRange($seg2_1.End, $search2.End)

# --- Paragraph 3: split into its 3 source runs ---
$search3 = $d.Paragraphs(3).Range
$seg3_0 = $search3
$null = $seg3_0.Find.Execute("CSS stands for Cascading style sheets. Here cascade means authors can define multiple styles in many different ways and these styles can be cascaded to many different levels where those levels may be of different documents or different elements in the same document. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg3_0.Font.Bold = 1
$seg3_0.Font.Bold = 0
$search3 = $d.Paragraphs(3).Range
$search3.SetRange($seg3_0.End, $search3.End)
$seg3_1 = $search3
$null = $seg3_1.Find.Execute("By 1996, CSS1 was finished while CSS working group started working on CSS2. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg3_1.Font.Bold = 1
$seg3_1.Font.Bold = 0
$search3 = $d.Paragraphs(3).Range
$search3.SetRange($seg3_1.End, $search3.End)
$seg3_2 = $search3
$null = $seg3_2.Find.Execute("As CSS started evolving, many different browsers started adopting it and now the implementation of CSS standards was being provided by almost all the browsers and CSS working group is progressing towards providing more mature effects from CSS.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg3_2.Font.Bold = 1
$seg3_2.Font.Bold = 0
$search3 = $d.Paragraphs(3).Range
$search3.SetRange($seg3_2.End, $search3.End)

# --- Paragraph 4: split into its 4 source runs ---
$search4 = $d.Paragraphs(4).Range
$seg4_0 = $search4
$null = $seg4_0.Find.Execute("Elements: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg4_0.Font.Bold = 1
$search4 = $d.Paragraphs(4).Range
$search4.SetRange($seg4_0.End, $search4.End)
$seg4_1 = $search4
$null = $seg4_1.Find.Execute("elements are the basis for document structure. In HTML, the most common elements are p, table, span, a, and div. Each of these elements ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg4_1.Font.Bold = 1
$seg4_1.Font.Bold = 0
$search4 = $d.Paragraphs(4).Range
$search4.SetRange($seg4_1.End, $search4.End)
$seg4_2 = $search4
$null = $seg4_2.Find.Execute("plays", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg4_2.Font.Bold = 1
$seg4_2.Font.Bold = 0
$search4 = $d.Paragraphs(4).Range
$search4.SetRange($seg4_2.End, $search4.End)
$seg4_3 = $search4
$null = $seg4_3.Find.Execute(" a role in the display of the document", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg4_3.Font.Bold = 1
$seg4_3.Font.Bold = 0
$search4 = $d.Paragraphs(4).Range
$search4.SetRange($seg4_3.End, $search4.End)

# --- Paragraph 5: split into its 2 source runs ---
$search5 = $d.Paragraphs(5).Range
$seg5_0 = $search5
$null = $seg5_0.Find.Execute("Replaced and Non-Replaced Elements: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg5_0.Font.Bold = 1
$search5 = $d.Paragraphs(5).Range
$search5.SetRange($seg5_0.End, $search5.End)
$seg5_1 = $search5
$null = $seg5_1.Find.Execute("Even though CSS depends on the elements of a document to decide how that element needs to be displayed on the screen, all the elements are not created equally. In CSS, elements take two forms", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg5_1.Font.Bold = 1
$seg5_1.Font.Bold = 0
$search5 = $d.Paragraphs(5).Range
$search5.SetRange($seg5_1.End, $search5.End)

# --- Paragraph 7: split into its 2 source runs ---
$search7 = $d.Paragraphs(7).Range
$seg7_0 = $search7
$null = $seg7_0.Find.Execute("Non-Replaced elements", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg7_0.Font.Bold = 1
$seg7_0.Font.Bold = 0
$search7 = $d.Paragraphs(7).Range
$search7.SetRange($seg7_0.End, $search7.End)
$seg7_1 = $search7
$null = $seg7_1.Find.Execute(": the majority of HTML elements are non-replaced elements. That means their content is presented by the user agent. Generally the user agent will be the browser. User agent presents the content inside a box that is generated by the element itself. For example, <span>hi</span> is a non-replaced element and the text “hi” will be displayed by the user agent. This is true for paragraphs, headings, table cells, lists and almost everything else in HTML.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg7_1.Font.Bold = 1
$seg7_1.Font.Bold = 0
$search7 = $d.Paragraphs(7).Range
$search7.SetRange($seg7_1.End, $search7.End)

# --- Paragraph 8: bold the whole (single-segment) paragraph ---
$p8 = $d.Paragraphs(8).Range
$p8.Font.Bold = 1

# --- Turn the 'Replaced elements' / 'Non-Replaced elements' paragraphs (6 & 7) into a numbered list ---
$d.Paragraphs(6).Range.set_Style("List Paragraph")
$d.Paragraphs(7).Range.set_Style("List Paragraph")
$listRange = $d.Range($d.Paragraphs(6).Range.Start, $d.Paragraphs(7).Range.End)
$listRange.ListFormat.ApplyNumberDefault()

# --- Match the generated 'List Paragraph' style's formatting to the target ---
$lpStyle = $d.Styles("List Paragraph")
$lpStyle.Priority = 34
$lpStyle.ParagraphFormat.LeftIndent = 36
$lpStyle.NoSpaceBetweenParagraphsOfSameStyle = 1

foreach ($p in $d.Paragraphs) {
    Write-Host "Para style=[" $p.Style.NameLocal "] text=[" $p.Range.Text "]"
}
